# Fruta / hortaliza, semanal
#
# Insert one new weekly price record as row 693 in the "Repollo" price
# history sheet. Excel shifts the existing rows 693:779 down to 694:780
# (carrying their formatting with them), which also grows the sheet's
# used range from A1:R779 to A1:R780 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the current row 693 (and everything below it) down by one row.
$ws.Rows(693).Insert()

# Populate the newly inserted row 693 with the new record.
$ws.Range("A693").Value = 4
$ws.Range("B693").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C693").Value = "Los Lagos"
$ws.Range("D693").Value = (Get-Date -Year 2023 -Month 7 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E693").Value = 10
$ws.Range("F693").Value = 100112006
$ws.Range("G693").Value = "Repollo"
$ws.Range("H693").Value = "Crespo record"
$ws.Range("I693").Value = "Primera"
$ws.Range("J693").Value = 500
$ws.Range("K693").Value = 1500
$ws.Range("L693").Value = 1700
$ws.Range("M693").Value = 1600
$ws.Range("N693").Value = "`$/unidad"
$ws.Range("O693").Value = "Región Metropolitana"
$ws.Range("P693").Value = 1600
$ws.Range("Q693").Value = 1
$ws.Range("R693").Value = "Hortaliza"
